$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 45727
$ws.Range("C1").Value = 45734
$ws.Range("D1").Value = 45741
$ws.Range("E1").Value = 45748
$ws.Range("F1").Value = 45755
$ws.Range("G1").Value = 45762
$ws.Range("H1").Value = 45769
$ws.Range("I1").Value = 45776
$ws.Range("J1").Value = 45783
$ws.Range("K1").Value = 45790
$ws.Range("L1").Value = 45797

$ws.Range("B2").Value = 48.1
$ws.Range("C2").Value = 49.2
$ws.Range("D2").Value = 50.5
$ws.Range("E2").Value = 51.7
$ws.Range("F2").Value = 53
$ws.Range("G2").Value = 54.1
$ws.Range("H2").Value = 55
$ws.Range("I2").Value = 55.5
$ws.Range("J2").Value = 55.7
$ws.Range("K2").Value = 55.6
$ws.Range("L2").Value = 55.1

$ws.Range("B3").Value = 59.8
$ws.Range("C3").Value = 61.5
$ws.Range("D3").Value = 62.8
$ws.Range("E3").Value = 63.7
$ws.Range("F3").Value = 64.1
$ws.Range("G3").Value = 64
$ws.Range("H3").Value = 63.2
$ws.Range("I3").Value = 61.8
$ws.Range("J3").Value = 59.9
$ws.Range("K3").Value = 57.6
$ws.Range("L3").Value = 55.2

$ws.Range("B4").Value = 40.5
$ws.Range("C4").Value = 41.8
$ws.Range("D4").Value = 43.6
$ws.Range("E4").Value = 45.9
$ws.Range("F4").Value = 48.7
$ws.Range("G4").Value = 51.7
$ws.Range("H4").Value = 54.8
$ws.Range("I4").Value = 57.6
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 61.9
$ws.Range("L4").Value = 63.1

$ws.Range("B5").Value = 37.6
$ws.Range("C5").Value = 39.5
$ws.Range("D5").Value = 41.8
$ws.Range("E5").Value = 44.5
$ws.Range("F5").Value = 47.4
$ws.Range("G5").Value = 50.3
$ws.Range("H5").Value = 53
$ws.Range("I5").Value = 55.3
$ws.Range("J5").Value = 57.1
$ws.Range("K5").Value = 58.3
$ws.Range("L5").Value = 58.8

$ws.Range("B6").Value = 54.6
$ws.Range("C6").Value = 54.4
$ws.Range("D6").Value = 54.3
$ws.Range("E6").Value = 54
$ws.Range("F6").Value = 53.7
$ws.Range("G6").Value = 53.3
$ws.Range("H6").Value = 52.6
$ws.Range("I6").Value = 51.7
$ws.Range("J6").Value = 50.6
$ws.Range("K6").Value = 49.3
$ws.Range("L6").Value = 48

$ws.Range("B7").Value = 49.6
$ws.Range("C7").Value = 54.1
$ws.Range("D7").Value = 58.4
$ws.Range("E7").Value = 62.3
$ws.Range("F7").Value = 65.7
$ws.Range("G7").Value = 68.6
$ws.Range("H7").Value = 71
$ws.Range("I7").Value = 72.9
$ws.Range("J7").Value = 74.3
$ws.Range("K7").Value = 75.1
$ws.Range("L7").Value = 75.5

$ws.Range("B8").Value = 45.2
$ws.Range("C8").Value = 46.3
$ws.Range("D8").Value = 47.6
$ws.Range("E8").Value = 49
$ws.Range("F8").Value = 50.4
$ws.Range("G8").Value = 51.6
$ws.Range("H8").Value = 52.6
$ws.Range("I8").Value = 53.3
$ws.Range("J8").Value = 53.4
$ws.Range("K8").Value = 53.2
$ws.Range("L8").Value = 52.6

$ws.Range("B9").Value = 59
$ws.Range("C9").Value = 62.2
$ws.Range("D9").Value = 65.1
$ws.Range("E9").Value = 67.5
$ws.Range("F9").Value = 69.5
$ws.Range("G9").Value = 70.9
$ws.Range("H9").Value = 71.8
$ws.Range("I9").Value = 72.1
$ws.Range("J9").Value = 72
$ws.Range("K9").Value = 71.4
$ws.Range("L9").Value = 70.4

$ws.Range("B10").Value = 31.8
$ws.Range("C10").Value = 33.1
$ws.Range("D10").Value = 34.9
$ws.Range("E10").Value = 37.3
$ws.Range("F10").Value = 40.2
$ws.Range("G10").Value = 43.4
$ws.Range("H10").Value = 46.6
$ws.Range("I10").Value = 49.5
$ws.Range("J10").Value = 52
$ws.Range("K10").Value = 53.9
$ws.Range("L10").Value = 55.1

$ws.Range("B11").Value = 52.8
$ws.Range("C11").Value = 51.2
$ws.Range("D11").Value = 49.7
$ws.Range("E11").Value = 48.1
$ws.Range("F11").Value = 46.7
$ws.Range("G11").Value = 45.2
$ws.Range("H11").Value = 43.7
$ws.Range("I11").Value = 42.1
$ws.Range("J11").Value = 40.4
$ws.Range("K11").Value = 38.8
$ws.Range("L11").Value = 37.3

$ws.Range("B12").Value = 55
$ws.Range("C12").Value = 57.1
$ws.Range("D12").Value = 59.2
$ws.Range("E12").Value = 61.4
$ws.Range("F12").Value = 63.4
$ws.Range("G12").Value = 65.3
$ws.Range("H12").Value = 66.9
$ws.Range("I12").Value = 68.1
$ws.Range("J12").Value = 68.8
$ws.Range("K12").Value = 69
$ws.Range("L12").Value = 68.7

$ws.Range("B13").Value = 65.4
$ws.Range("C13").Value = 65.6
$ws.Range("D13").Value = 65.5
$ws.Range("E13").Value = 65
$ws.Range("F13").Value = 64.1
$ws.Range("G13").Value = 62.7
$ws.Range("H13").Value = 60.7
$ws.Range("I13").Value = 58.2
$ws.Range("J13").Value = 55.2
$ws.Range("K13").Value = 52
$ws.Range("L13").Value = 48.6

$ws.Range("B14").Value = 51
$ws.Range("C14").Value = 49.2
$ws.Range("D14").Value = 46.8
$ws.Range("E14").Value = 43.9
$ws.Range("F14").Value = 40.4
$ws.Range("G14").Value = 36.4
$ws.Range("H14").Value = 32
$ws.Range("I14").Value = 27.5
$ws.Range("J14").Value = 23.1
$ws.Range("K14").Value = 19.2
$ws.Range("L14").Value = 15.9

$ws.Range("B15").Value = 71
$ws.Range("C15").Value = 73.7
$ws.Range("D15").Value = 75.8
$ws.Range("E15").Value = 77.5
$ws.Range("F15").Value = 78.5
$ws.Range("G15").Value = 79
$ws.Range("H15").Value = 79.1
$ws.Range("I15").Value = 78.7
$ws.Range("J15").Value = 78.1
$ws.Range("K15").Value = 77.1
$ws.Range("L15").Value = 76

$ws.Range("B16").Value = 32.3
$ws.Range("C16").Value = 33.7
$ws.Range("D16").Value = 35.5
$ws.Range("E16").Value = 37.7
$ws.Range("F16").Value = 40.4
$ws.Range("G16").Value = 43.2
$ws.Range("H16").Value = 46.1
$ws.Range("I16").Value = 48.7
$ws.Range("J16").Value = 50.9
$ws.Range("K16").Value = 52.7
$ws.Range("L16").Value = 53.9

$ws.Range("B17").Value = 33.2
$ws.Range("C17").Value = 33.3
$ws.Range("D17").Value = 33.8
$ws.Range("E17").Value = 34.5
$ws.Range("F17").Value = 35.6
$ws.Range("G17").Value = 36.6
$ws.Range("H17").Value = 37.6
$ws.Range("I17").Value = 38.3
$ws.Range("J17").Value = 38.8
$ws.Range("K17").Value = 39
$ws.Range("L17").Value = 39.1

$ws.Range("B18").Value = 40.8
$ws.Range("C18").Value = 41.6
$ws.Range("D18").Value = 43.1
$ws.Range("E18").Value = 45.2
$ws.Range("F18").Value = 48.1
$ws.Range("G18").Value = 51.5
$ws.Range("H18").Value = 55.1
$ws.Range("I18").Value = 58.5
$ws.Range("J18").Value = 61.6
$ws.Range("K18").Value = 64.1
$ws.Range("L18").Value = 65.8

$ws.Range("B19").Value = 38
$ws.Range("C19").Value = 38.9
$ws.Range("D19").Value = 40.1
$ws.Range("E19").Value = 41.7
$ws.Range("F19").Value = 43.7
$ws.Range("G19").Value = 45.8
$ws.Range("H19").Value = 47.9
$ws.Range("I19").Value = 49.9
$ws.Range("J19").Value = 51.5
$ws.Range("K19").Value = 52.6
$ws.Range("L19").Value = 53.1
